# This script updates the LR-pair TPM output values in the active sheet
# to reflect a re-run of the NATMI pipeline with new TPM data.
# Columns (1-indexed):
#  A=1 Sending cluster, B=2 Ligand symbol, C=3 Receptor symbol, D=4 Target cluster
#  E=5 Ligand-expressing cells, F=6 Ligand detection rate
#  G=7 Ligand average expression value, H=8 Ligand total expression value
#  I=9 Ligand derived specificity (avg), J=10 Ligand derived specificity (total)
#  K=11 Receptor-expressing cells, L=12 Receptor detection rate
#  M=13 Receptor average expression value, N=14 Receptor total expression value
#  O=15 Receptor derived specificity (avg), P=16 Receptor derived specificity (total)
#  Q=17 Edge average expression weight, R=18 Edge total expression weight
#  S=19 Edge average expression derived specificity, T=20 Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ E=3; F=1; G=15.58971733333333; H=46.769152; I=0.58284954614551; J=0.5828495461455101; K=3; L=1; M=5.988024; N=17.964072; O=0.4288861341243614; P=0.4288861341243613; Q=93.351601545216; R=840.164413906944; S=0.2499760886224864; T=0.2499760886224864 }
    3  = @{ E=3; F=1; G=15.58971733333333; H=46.769152; I=0.58284954614551; J=0.5828495461455101;                                                                             O=0.5018175899002161; P=0.5018175899002161; Q=109.2259039718969; R=983.0331357470719; S=0.2924841545211746; T=0.2924841545211747 }
    4  = @{ E=3; F=1; G=15.58971733333333; H=46.769152; I=0.58284954614551; J=0.5828495461455101;                        M=0.9675009999999999; N=2.902503; O=0.06929627597542257; P=0.06929627597542255; Q=15.08306710971733; R=135.747603987456; S=0.04038930300184905; T=0.04038930300184905 }
    5  = @{                                                                 I=0.3268220285680467; J=0.3268220285680467; K=3; L=1; M=5.988024; N=17.964072; O=0.4288861341243614; P=0.4288861341243613; Q=52.345172075448; R=471.106548679032; S=0.1401694363792311; T=0.1401694363792311 }
    6  = @{                                                                 I=0.3268220285680467; J=0.3268220285680467;                                                                             O=0.5018175899002161; P=0.5018175899002161;                                                                             S=0.1640050427023168; T=0.1640050427023168 }
    7  = @{                                                                 I=0.3268220285680467; J=0.3268220285680467;                                                    M=0.9675009999999999; N=2.902503; O=0.06929627597542257; P=0.06929627597542255; Q=8.457548989143666; R=76.11794090229299; S=0.0226475494864988; T=0.0226475494864988 }
    8  = @{                                  G=2.416051666666667; H=7.248155; I=0.09032842528644328; J=0.09032842528644329; K=3; L=1; M=5.988024; N=17.964072; O=0.4288861341243614; P=0.4288861341243613; Q=14.46737536524; R=130.20637828716; S=0.03874060912264387; T=0.03874060912264387 }
    9  = @{                                  G=2.416051666666667; H=7.248155; I=0.09032842528644328; J=0.09032842528644329;                                                                             O=0.5018175899002161; P=0.5018175899002161; Q=16.92753124973111; R=152.34778124758; S=0.0453283926767247; T=0.04532839267672471 }
    10 = @{                                  G=2.416051666666667; H=7.248155; I=0.09032842528644328; J=0.09032842528644329;                                                    M=0.9675009999999999; N=2.902503; O=0.06929627597542257; P=0.06929627597542255; Q=2.337532403551667; R=21.037791631965; S=0.006259423487074711; T=0.006259423487074711 }
}

$colIndex = @{ E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($colLetter in $rowData.Keys) {
        $colNum = $colIndex[$colLetter]
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$colLetter]
    }
}
